# Apply updated "dSF" (column F) values for the rows identified in the diff.
# Mapping is keyed by worksheet row number -> new value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    3  = -4
    5  = 1
    8  = -1
    13 = 0
    15 = 1
    17 = -1
    25 = 1
    34 = -1
    36 = 5
    42 = 4
    51 = -1
    52 = 2
    56 = -1
    57 = 2
    59 = 3
    65 = 2
    66 = -3
    69 = -2
    70 = -2
    72 = -3
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
